$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text formatting so numeric-looking values
# (e.g. "168.10", "1.18") are not auto-converted to numbers, which would
# lose formatting (trailing zeros) and change the stored type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.823.65'
$ws.Range('E2').Value = '  +3.83%  '
$ws.Range('D3').Value = '2.276.81'
$ws.Range('E3').Value = '  +4.27%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '251.15'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = '0.636'
$ws.Range('E6').Value = '  +3.98%  '
$ws.Range('D7').Value = '71.70'
$ws.Range('E7').Value = '  +7.72%  '
$ws.Range('D9').Value = '0.647'
$ws.Range('E9').Value = '  +12.88%  '
$ws.Range('D10').Value = '38.45'
$ws.Range('E10').Value = '  +5.40%  '
$ws.Range('D11').Value = '59.81'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').Value = '0.0971'
$ws.Range('E12').Value = '  +4.64%  '
$ws.Range('D13').Value = '7.36'
$ws.Range('E13').Value = '  +6.47%  '
$ws.Range('D14').Value = '0.105'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('D15').Value = '2.617.42'
$ws.Range('E15').Value = '  +4.43%  '
$ws.Range('D16').Value = '14.89'
$ws.Range('E16').Value = '  +3.47%  '
$ws.Range('D17').Value = '0.885'
$ws.Range('E17').Value = '  +3.59%  '
$ws.Range('D18').Value = '2.275.99'
$ws.Range('E18').Value = '  +4.85%  '
$ws.Range('D19').Value = '42.780.31'
$ws.Range('E19').Value = '  +3.92%  '
$ws.Range('D20').Value = '0.0000101'
$ws.Range('E20').Value = '  +7.14%  '
$ws.Range('D21').Value = '6.30'
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('D22').Value = '73.43'
$ws.Range('E22').Value = '  +2.40%  '
$ws.Range('D23').Value = '233.93'
$ws.Range('E23').Value = '  +1.45%  '
$ws.Range('D24').Value = '2.13'
$ws.Range('E24').Value = '  +4.64%  '
$ws.Range('D25').Value = '4.03'
$ws.Range('E25').Value = '  +7.05%  '
$ws.Range('D26').Value = '11.48'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = '2.44'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '3.67'
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('E30').Value = '  +5.21%  '
$ws.Range('D31').Value = '168.10'
$ws.Range('D32').Value = '20.97'
$ws.Range('E32').Value = '  +3.44%  '
$ws.Range('D33').Value = '6.46'
$ws.Range('E33').Value = '  +11.88%  '
$ws.Range('E34').Value = '  +4.03%  '
$ws.Range('D35').Value = '31.49'
$ws.Range('E35').Value = '  +27.73%  '
$ws.Range('D36').Value = '0.0792'
$ws.Range('E36').Value = '  +7.45%  '
$ws.Range('E37').Value = '  +3.55%  '
$ws.Range('D38').Value = '4.52'
$ws.Range('E38').Value = '  +13.30%  '
$ws.Range('D39').Value = '4.76'
$ws.Range('E39').Value = '  +4.90%  '
$ws.Range('D40').Value = '0.0314'
$ws.Range('E40').Value = '  +2.15%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '2.32'
$ws.Range('E41').Value = '  +4.59%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').Value = '13.21'
$ws.Range('E42').Value = '  +15.19%  '
$ws.Range('D43').Value = '5.83'
$ws.Range('E43').Value = '  +5.87%  '
$ws.Range('D44').Value = '0.209'
$ws.Range('E44').Value = '  +9.52%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = '5.04'
$ws.Range('E45').Value = '  -6.49%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '9.17'
$ws.Range('E46').Value = '  +7.55%  '
$ws.Range('B47').Value = 'MultiversX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D47').Value = '61.94'
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').Value = '1.18'
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('D51').Value = '1.19'
$ws.Range('E51').Value = '  +3.53%  '
